$wb = $excel.ActiveWorkbook

# --- "Repayment schedule" sheet (sheet5): insert a new blank column before N ---
# This shifts the old N/O/P ("Late" header + data, blank spacer, "Outstanding"
# header + data) one column to the right, landing on O/P/Q, and leaves a new
# blank (but styled) column N in their place - matching the column-N insert
# visible in the diff (headers + all 14 data rows for columns N..P shift to
# O..Q, dimension grows from A1:P14 to A1:Q14, row spans 1:16 -> 1:17).
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").Insert()

# The insert above pushes column C's bestFit computation slightly (and the
# newly created column D now holds content too) - nudge both to their
# recalculated best-fit widths.
$wsSchedule.Columns.Item(3).ColumnWidth = 9.333333333333332
$wsSchedule.Columns.Item(4).ColumnWidth = 9

# "Repayment schedule" becomes the active/selected sheet (tabSelected moves
# here from "Edit Repayment Schedule"), with the cursor resting on E13.
$wsSchedule.Activate()
[void]$wsSchedule.Range("E13").Select()
